# Applies the "Updated cryptos list" diff to Sheet1 of the workbook.
# For each changed cell we simply write the new literal text. Columns D/E
# and B/C hold plain display strings (prices, percent deltas, coin
# names/links) that must stay text, exactly as authored upstream
# (t="inlineStr" in the source XML) - e.g. "1.003" or "225.61" must NOT
# turn into the numbers 1.003 / 225.61. Excel's Range.Value setter auto-
# detects numeric-looking strings and coerces them, so for any new value
# that parses as a plain number we temporarily force the cell to Text
# format, assign the value, then clear the formatting again so the cell's
# style index is left exactly as it was (no stray numFmt/style changes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $value) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# Row 2
$ws.Range("D2").Value = "27.527.15"
$ws.Range("E2").Value = "  +5.58%  "
# Row 3
$ws.Range("D3").Value = "1.722.08"
$ws.Range("E3").Value = "  +4.33%  "
# Row 4
Set-TextCell "D4" "1.003"
$ws.Range("E4").Value = "  +0.10%  "
# Row 5
Set-TextCell "D5" "225.61"
$ws.Range("E5").Value = "  +3.45%  "
# Row 6
Set-TextCell "D6" "0.5367"
$ws.Range("E6").Value = "  +3.01%  "
# Row 7
Set-TextCell "D7" "1.003"
$ws.Range("E7").Value = "  +0.05%  "
# Row 8
Set-TextCell "D8" "0.2663"
$ws.Range("E8").Value = "  +1.13%  "
# Row 9
Set-TextCell "D9" "0.06600"
$ws.Range("E9").Value = "  +4.37%  "
# Row 10
Set-TextCell "D10" "21.67"
$ws.Range("E10").Value = "  +6.40%  "
# Row 11
Set-TextCell "D11" "0.07711"
$ws.Range("E11").Value = "  +0.76%  "
# Row 12
Set-TextCell "D12" "4.617"
$ws.Range("E12").Value = "  +0.63%  "
# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.717.46"
$ws.Range("E13").Value = "  +4.24%  "
# Row 14
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "1.960.46"
$ws.Range("E14").Value = "  +4.43%  "
# Row 15
Set-TextCell "D15" "0.5834"
$ws.Range("E15").Value = "  +4.42%  "
# Row 16
$ws.Range("D16").Value = "0.0₅8301"
$ws.Range("E16").Value = "  +2.04%  "
# Row 17
Set-TextCell "D17" "67.88"
$ws.Range("E17").Value = "  +4.10%  "
# Row 18
$ws.Range("D18").Value = "27.536.76"
$ws.Range("E18").Value = "  +5.69%  "
# Row 19
Set-TextCell "D19" "220.04"
$ws.Range("E19").Value = "  +15.15%  "
# Row 20
$ws.Range("E20").Value = "  +0.08%  "
# Row 21
Set-TextCell "D21" "4.725"
$ws.Range("E21").Value = "  +2.22%  "
# Row 22
Set-TextCell "D22" "10.64"
$ws.Range("E22").Value = "  +1.65%  "
# Row 23
Set-TextCell "D23" "6.077"
$ws.Range("E23").Value = "  +2.93%  "
# Row 24
Set-TextCell "D24" "1.004"
$ws.Range("E24").Value = "  +0.10%  "
# Row 25
Set-TextCell "D25" "148.14"
$ws.Range("E25").Value = "  +2.79%  "
# Row 26
$ws.Range("E26").Value = "  +14.86%  "
# Row 27
Set-TextCell "D27" "0.1235"
$ws.Range("E27").Value = "  +4.13%  "
# Row 28
Set-TextCell "D28" "7.398"
$ws.Range("E28").Value = "  +2.76%  "
# Row 29
Set-TextCell "D29" "16.62"
$ws.Range("E29").Value = "  +4.69%  "
# Row 30
Set-TextCell "D30" "0.05575"
$ws.Range("E30").Value = "  +2.63%  "
# Row 31
Set-TextCell "D31" "1.301"
$ws.Range("E31").Value = "  +2.73%  "
# Row 32
$ws.Range("E32").Value = "  +3.62%  "
# Row 33
Set-TextCell "D33" "3.443"
$ws.Range("E33").Value = "  +2.93%  "
# Row 34
Set-TextCell "D34" "1.665"
$ws.Range("E34").Value = "  +7.13%  "
# Row 35
Set-TextCell "D35" "2.839"
$ws.Range("E35").Value = "  +2.06%  "
# Row 36
Set-TextCell "D36" "0.9620"
# Row 37
Set-TextCell "D37" "2.429"
$ws.Range("E37").Value = "  +0.18%  "
# Row 38
Set-TextCell "D38" "0.5966"
$ws.Range("E38").Value = "  +6.00%  "
# Row 39
$ws.Range("E39").Value = "  +4.33%  "
# Row 40
$ws.Range("E40").Value = "  +1.15%  "
# Row 41
Set-TextCell "D41" "0.8553"
$ws.Range("E41").Value = "  +3.56%  "
# Row 42
$ws.Range("D42").Value = "1.053.36"
$ws.Range("E42").Value = "  +2.57%  "
# Row 43
Set-TextCell "D43" "1.003"
$ws.Range("E43").Value = "  +0.07%  "
# Row 44
Set-TextCell "D44" "101.28"
$ws.Range("E44").Value = "  +0.20%  "
# Row 45
$ws.Range("D45").Value = "1.866.55"
$ws.Range("E45").Value = "  +4.54%  "
# Row 46
$ws.Range("E46").Value = "  +4.07%  "
# Row 47
Set-TextCell "D47" "59.12"
$ws.Range("E47").Value = "  +3.21%  "
# Row 48
Set-TextCell "D48" "8.220"
$ws.Range("E48").Value = "  +3.54%  "
# Row 49
Set-TextCell "D49" "0.4434"
$ws.Range("E49").Value = "  +2.41%  "
# Row 50
Set-TextCell "D50" "1.000"
$ws.Range("E50").Value = "  +0.01%  "
# Row 51
Set-TextCell "D51" "0.05244"
$ws.Range("E51").Value = "  +1.95%  "
